# "change valeu menu tengah"
# Update the two "menu tengah" (middle menu) lookup sheets so that the
# "News" (News/trending) row becomes "Drama" and the "Radio+" (Radio+/radio)
# row becomes "Comedy". Also rename the "Menu_Pilar_Not" sheet to
# "Menu_Pilar_Tak_Tampil" and nudge a few sheet selections/active-tab state
# to match the author's final click position when they saved the file.

$wb = $excel.ActiveWorkbook

# --- Homepage_Menu_Tengah: B3 News -> Drama, B4 Radio+ -> Comedy ---
$wsTengah = $wb.Worksheets.Item("Homepage_Menu_Tengah")
$wsTengah.Range("B3").Value = "Drama"
$wsTengah.Range("B4").Value = "Comedy"
$wsTengah.Range("B2:B4").Select()

# --- Homepage_Menu_Tengah_Direct: same value swap in both the label (B)
#     and slug (C) columns ---
$wsTengahDirect = $wb.Worksheets.Item("Homepage_Menu_Tengah_Direct")
$wsTengahDirect.Range("B3").Value = "Drama"
$wsTengahDirect.Range("C3").Value = "Drama"
$wsTengahDirect.Range("B4").Value = "Comedy"
$wsTengahDirect.Range("C4").Value = "Comedy"

# --- Rename Menu_Pilar_Not -> Menu_Pilar_Tak_Tampil and move its cursor ---
$wsPilarNot = $wb.Worksheets.Item("Menu_Pilar_Not")
$wsPilarNot.Name = "Menu_Pilar_Tak_Tampil"
$wsPilarNot.Range("H11").Select()

# --- Menu_Pilar_Tampil: cursor moved from K14 to K6 ---
$wsPilarTampil = $wb.Worksheets.Item("Menu_Pilar_Tampil")
$wsPilarTampil.Range("K6").Select()

# --- Finish on Homepage_Menu_Tengah_Direct, which becomes the active tab ---
$wsTengahDirect.Activate()
$wsTengahDirect.Range("J9").Select()
